$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D=44161; J=270; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=3; D=44468; J=300; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=4; D=44243; J=250; K=1200; L=1300; M=1250; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=625; Q=2},
    @{Row=5; D=44438; J=300; K=950; L=1000; M=975; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=488; Q=2},
    @{Row=6; D=44390; J=250; K=2400; L=2500; M=2450; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=1225; Q=2},
    @{Row=7; D=44257; J=500; K=1400; L=1500; M=1450; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=725; Q=2},
    @{Row=8; D=44266; J=300; K=1700; L=1800; M=1750; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=875; Q=2},
    @{Row=9; D=44385; J=300; K=2400; L=2500; M=2450; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=1225; Q=2},
    @{Row=10; D=44403; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=11; D=44363; J=250; K=2500; L=2800; M=2650; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=1325; Q=2},
    @{Row=12; D=44253; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=13; D=44435; J=300; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=14; D=44302; J=300; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=15; D=44181; J=200; K=1000; L=1200; M=1100; N='$/atado'; O='Región de Arica y Parinacota'; P=1100; Q=1},
    @{Row=16; D=44291; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=17; D=44616; J=270; K=1300; L=1500; M=1400; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=700; Q=2},
    @{Row=18; D=44540; J=300; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=19; D=44229; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=20; D=44202; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=21; D=44525; J=300; K=1400; L=1500; M=1450; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=725; Q=2},
    @{Row=22; D=44447; J=300; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=23; D=44172; J=200; K=1300; L=1500; M=1400; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=700; Q=2},
    @{Row=24; D=44601; J=270; K=2200; L=2500; M=2350; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=1175; Q=2},
    @{Row=25; D=44365; J=200; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2},
    @{Row=26; D=44427; J=250; K=1300; L=1500; M=1400; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=700; Q=2},
    @{Row=27; D=44572; J=300; K=1400; L=1500; M=1450; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=725; Q=2},
    @{Row=28; D=44544; J=250; K=900; L=1000; M=950; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=475; Q=2},
    @{Row=29; D=44392; J=250; K=1800; L=2000; M=1900; N='$/atado 1,5 a 2 kilos'; O='Región de Arica y Parinacota'; P=950; Q=2}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 10).Value = $item.J
    $ws.Cells.Item($item.Row, 11).Value = $item.K
    $ws.Cells.Item($item.Row, 12).Value = $item.L
    $ws.Cells.Item($item.Row, 13).Value = $item.M
    $ws.Cells.Item($item.Row, 14).Value = $item.N
    $ws.Cells.Item($item.Row, 15).Value = $item.O
    $ws.Cells.Item($item.Row, 16).Value = $item.P
    $ws.Cells.Item($item.Row, 17).Value = $item.Q
}